$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: "Butir" header + blank, bordered body cells for rows 2-29.
# Clone formatting from the existing last header/body column (J) so the
# existing style indices (s="1" header, s="2" body) are reused instead of
# new styles/fills being synthesized.
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "Butir"

$ws.Range("J2:J29").Copy($ws.Range("K2:K29"))
$ws.Range("K2:K29").ClearContents()

# Widen the new column.
$ws.Columns.Item(11).ColumnWidth = 14.83

# Match the selection left behind in the authored workbook.
$ws.Range("O20").Select()
